# Update gh-pages output figures (F column "浏览量/人气" counters) for
# several rows across the "展览", "演出" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")     # sheet1
$wsShow    = $wb.Worksheets.Item("演出")     # sheet2
$wsAll     = $wb.Worksheets.Item("全部类型") # sheet4

# -- 展览 (sheet1) --
$wsExhibit.Range("F3").Value  = 393
$wsExhibit.Range("F6").Value  = 133
$wsExhibit.Range("F7").Value  = 875
$wsExhibit.Range("F8").Value  = 681
$wsExhibit.Range("F9").Value  = 139
$wsExhibit.Range("F12").Value = 744
$wsExhibit.Range("F15").Value = 474
$wsExhibit.Range("F19").Value = 2735
$wsExhibit.Range("F20").Value = 1207
$wsExhibit.Range("F21").Value = 624
$wsExhibit.Range("F24").Value = 50
$wsExhibit.Range("F25").Value = 933
$wsExhibit.Range("F27").Value = 1256

# -- 演出 (sheet2) --
$wsShow.Range("F3").Value = 498

# -- 全部类型 (sheet4) --
$wsAll.Range("F5").Value  = 393
$wsAll.Range("F7").Value  = 498
$wsAll.Range("F8").Value  = 498
$wsAll.Range("F13").Value = 133
$wsAll.Range("F14").Value = 875
$wsAll.Range("F15").Value = 681
$wsAll.Range("F16").Value = 139
$wsAll.Range("F24").Value = 744
$wsAll.Range("F27").Value = 474
$wsAll.Range("F31").Value = 2735
$wsAll.Range("F32").Value = 1207
$wsAll.Range("F33").Value = 624
$wsAll.Range("F36").Value = 50
$wsAll.Range("F38").Value = 933
$wsAll.Range("F40").Value = 1256
